$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in column A (ID) for both new rows first, matching the order the
# author typed the new data (IDs entered before the rest of each row)
$ws.Range("A22").Value = "IA_0021"
$ws.Range("A23").Value = "IA_0022"

# Row 22 - Ginny playing in the snow
$ws.Range("B22").Value = "IMG_5877.jpg"
$ws.Range("D22").Value = "Nora Ammann"
$ws.Range("E22").Value = "CC BY 4.0"
$ws.Range("F22").Value = "data/Multimedia_Data/Image_Animal/"
$ws.Range("G22").Value = "Ginny playing in the snow"
$ws.Range("H22").Value = "A_001"
$ws.Range("I22").Value = 18

# Row 23 - Ginny sleeping
$ws.Range("B23").Value = "IMG_7498.jpg"
$ws.Range("D23").Value = "Nora Ammann"
$ws.Range("E23").Value = "CC BY 4.0"
$ws.Range("F23").Value = "data/Multimedia_Data/Image_Animal/"
$ws.Range("G23").Value = "Ginny sleeping"
$ws.Range("H23").Value = "A_001"
$ws.Range("I23").Value = 19

$ws.Range("H27").Select() | Out-Null
